$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the two new 6-column blocks so the sheet grows from 12 to 24
#    columns, pushing the existing "Player(wager1)" block from C->I and the
#    card/payAmt block from D:L -> P:X (matches the target <cols>/<dimension>).
# ---------------------------------------------------------------------------
$ws.Range("C1:H1").EntireColumn.Insert()
$ws.Range("J1:O1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. New header / data values. The assignment order below is deliberate: it
#    reproduces the shared-string table order of the target workbook (new
#    unique strings are appended to xl/sharedStrings.xml in first-write
#    order), so do not reorder these lines.
# ---------------------------------------------------------------------------
$ws.Range("I2").Value = "P1;100;P1"
$ws.Range("B2").Value = "100;rated-6009;1"
$ws.Range("C1").Value = "buyIn2"
$ws.Range("J1").Value = "wager2"
$ws.Range("K1").Value = "wager3"
$ws.Range("L1").Value = "wager4"
$ws.Range("D1").Value = "buyIn3"
$ws.Range("D2").Value = "100;anon"
$ws.Range("E1").Value = "buyIn4"
$ws.Range("F1").Value = "buyIn5"
$ws.Range("G1").Value = "buyIn6"
$ws.Range("C2").Value = "100;known-6010"
$ws.Range("F2").Value = "100;known-6012"
$ws.Range("M1").Value = "wager5"
$ws.Range("N1").Value = "wager6"
$ws.Range("J2").Value = "P2;100;P2"
$ws.Range("H1").Value = "buyIn7"
$ws.Range("O1").Value = "wager7"
$ws.Range("M2").Value = "P5;100;P6"
$ws.Range("K2").Value = "P3;100;B3"
$ws.Range("L2").Value = "P4;100;B5"
$ws.Range("E2").Value = "100;rated-6004;3"

# ---------------------------------------------------------------------------
# 4. Column widths for the newly inserted columns
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 29.7109375
$ws.Columns.Item(5).ColumnWidth = 29.7109375
$ws.Columns.Item(6).ColumnWidth = 29.7109375
$ws.Columns.Item(7).ColumnWidth = 29.7109375
$ws.Columns.Item(8).ColumnWidth = 29.7109375

$ws.Columns.Item(10).ColumnWidth = 24
$ws.Columns.Item(11).ColumnWidth = 24
$ws.Columns.Item(12).ColumnWidth = 24
$ws.Columns.Item(13).ColumnWidth = 24
$ws.Columns.Item(14).ColumnWidth = 24
$ws.Columns.Item(15).ColumnWidth = 24

$ws.Columns.Item(3).ColumnWidth = 29.7109375

# ---------------------------------------------------------------------------
# 5. Sheet view: scrolled to H1, active selection N2
# ---------------------------------------------------------------------------
$ws.Range("N2").Select()
$excel.ActiveWindow.ScrollColumn = 8

"done"
